$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H ("Comment") gets new explanatory text for most question types ---
# NB: new shared-string entries are appended in the order the cells are
# written below, so this order matters (it must match the order the
# strings appear in xl/sharedStrings.xml of the target workbook).

# Row 4 (Character type)
$ws.Range("H4").Value = "A single character is expected and only the first character is evaluated. Capitalization is ignored."

# Row 2 (Text type)
$ws.Range("H2").Value = "Leading and trailing white spaces are trimmed in the submission and capitalization is ignored, but otherwise matching is strict."

# Row 3 (Logical type) - rich text with a couple of runs in a smaller "Arial Unicode MS" font
$h3text = "Expected inputs are True/False, Yes/No, 1/0. Only the first character is evaluated (i.e.  T/Y/1 or F/N/0). Capitalization is ignored. Any non-match is considered a mistake."
$ws.Range("H3").Value = $h3text

$run2 = $ws.Range("H3").Characters(21, 19)
$run2.Font.Name = "Arial Unicode MS"
$run2.Font.Size = 10

$run3 = $ws.Range("H3").Characters(40, 1)
$run3.Font.Name = "Calibri"
$run3.Font.Size = 11

$run4 = $ws.Range("H3").Characters(41, 3)
$run4.Font.Name = "Arial Unicode MS"
$run4.Font.Size = 10

$run5 = $ws.Range("H3").Characters(44, 129)
$run5.Font.Name = "Calibri"
$run5.Font.Size = 11

# Row 5 (Integer type)
$ws.Range("H5").Value = "Expect an integer number. Anything non-integer or non-matching is considered incorrect."

# Row 6 (Fuzzy_Integer type) - reword the existing comment
$ws.Range("H6").Value = "This  allows wrong rounding of the answer to e.g. 41"

# Row 8 (Rounded_Numeric type)
$ws.Range("H8").Value = "A numeric value rounded to some specified significant digits. Those digits need to match the provided answer."

# --- Page setup: force portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to A8 ---
$ws.Range("A8").Select()
